$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ============================================================
# 1) Update "Weekly Timesheet" sheet data
# ============================================================

# Row 2: 2026-01-12
$ws1.Range("B2").Value = "Hall"
$ws1.Range("C2").Value = 9
$ws1.Range("F2").Value = 900

# Row 3: 2026-01-13
$ws1.Range("B3").Value = "Bryan"
$ws1.Range("C3").Value = 8
$ws1.Range("E3").Value = 90
$ws1.Range("F3").Value = 720

# Row 4: 2026-01-14
$ws1.Range("B4").Value = "McGill"
$ws1.Range("C4").Value = 7
$ws1.Range("E4").Value = 90
$ws1.Range("F4").Value = 630

# Row 5: 2026-01-15
$ws1.Range("B5").Value = "Hall"
$ws1.Range("C5").Value = 8
$ws1.Range("F5").Value = 800

# Row 6: 2026-01-16
$ws1.Range("B6").Value = "Bryan"
$ws1.Range("C6").Value = 10
$ws1.Range("E6").Value = 90
$ws1.Range("F6").Value = 900

# Row 8: SUBTOTAL
$ws1.Range("C8").Value = 42
$ws1.Range("D8").Value = "Reg: 42 / OT: 0"
$ws1.Range("F8").Value = 3950

# ============================================================
# 2) Add the new "Jason Schema" sheet right after "Weekly Timesheet"
# ============================================================

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Jason Schema"

# Column widths (roughly matching target layout)
$ws2.Columns.Item(1).ColumnWidth = 20
$ws2.Columns.Item(2).ColumnWidth = 18
$ws2.Columns.Item(3).ColumnWidth = 12
$ws2.Columns.Item(4).ColumnWidth = 25
$ws2.Columns.Item(5).ColumnWidth = 8
$ws2.Columns.Item(6).ColumnWidth = 10
$ws2.Columns.Item(7).ColumnWidth = 12
$ws2.Columns.Item(8).ColumnWidth = 10
$ws2.Columns.Item(9).ColumnWidth = 30

# Header row (bold)
$headers = @("Employee", "Employee ID", "Date", "Client", "Hours", "Rate", "Total", "Type", "Notes")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws2.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}
$ws2.Range("F1").NumberFormat = """$""#,##0.00"
$ws2.Range("G1").NumberFormat = """$""#,##0.00"

# Data rows: Employee, Employee ID, Date, Client, Hours, Rate, Total, Type
$rows = @(
    @("Chris Z", "emp_JcYsCv7rJ7fyha2O", "2026-01-12", "Hall",   9,  100, 900, "Regular"),
    @("Chris Z", "emp_JcYsCv7rJ7fyha2O", "2026-01-13", "Bryan",  8,  90,  720, "Regular"),
    @("Chris Z", "emp_JcYsCv7rJ7fyha2O", "2026-01-14", "McGill", 7,  90,  630, "Regular"),
    @("Chris Z", "emp_JcYsCv7rJ7fyha2O", "2026-01-15", "Hall",   8,  100, 800, "Regular"),
    @("Chris Z", "emp_JcYsCv7rJ7fyha2O", "2026-01-16", "Bryan",  10, 90,  900, "Regular")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $excelRow = $r + 2

    $ws2.Cells.Item($excelRow, 1).Value = $rowData[0]
    $ws2.Cells.Item($excelRow, 2).Value = $rowData[1]

    # Force text storage for the date column so it matches the source
    # sheet's convention of literal "YYYY-MM-DD" strings (not date serials)
    $dateCell = $ws2.Cells.Item($excelRow, 3)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[2]

    $ws2.Cells.Item($excelRow, 4).Value = $rowData[3]
    $ws2.Cells.Item($excelRow, 5).Value = $rowData[4]

    $rateCell = $ws2.Cells.Item($excelRow, 6)
    $rateCell.Value = $rowData[5]
    $rateCell.NumberFormat = """$""#,##0.00"

    $totalCell = $ws2.Cells.Item($excelRow, 7)
    $totalCell.Value = $rowData[6]
    $totalCell.NumberFormat = """$""#,##0.00"

    $ws2.Cells.Item($excelRow, 8).Value = $rowData[7]
}

Write-Host "Edit complete"
